# Regenerate merged AHB files
# 1) Rename the header cells (A1:U1) from the "_old"/"_new" diff-suffix
#    convention to the "_FV2404"/"_FV2410" convention (the "diff" column
#    header stays as-is).
# 2) Turn the used range A1:U72 into an Excel Table ("Table1") with no
#    explicit table style so banding/AutoFilter come from the table
#    defaults.
# 3) Freeze the header row (row 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410",
  "Segmentgruppe_FV2410",
  "Segment_FV2410",
  "Datenelement_FV2410",
  "Segment ID_FV2410",
  "Code_FV2410",
  "Qualifier_FV2410",
  "Beschreibung_FV2410",
  "Bedingungsausdruck_FV2410",
  "Bedingung_FV2410"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Convert A1:U72 into a real Table (ListObject) so the new xl/tables/table1.xml
# part + <tableParts> reference show up, matching a fresh "Insert Table".
$range = $ws.Range("A1:U72")
$tbl = $ws.ListObjects.Add(1, $range, $null, 1, $null)
$tbl.Name = "Table1"

# Freeze the header row: select the first row below the header, then freeze.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
